# Backlog.xlsx update: refactor(gameloop) — update sprint hours / status on the
# ARCHIVE (Backlog) sheet and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHIVE")

# --- Update logged HOURS (column G) for the current sprint rows ---
$ws.Range("G3").Value  = 0.5
$ws.Range("G4").Value  = 0.5
$ws.Range("G5").Value  = 0.5
$ws.Range("G6").Value  = 0.5
$ws.Range("G7").Value  = 0.5
$ws.Range("G8").Value  = 0.5
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 0.5
$ws.Range("G12").Value = 0.5
$ws.Range("G13").Value = 1
$ws.Range("G14").Value = 0.5
$ws.Range("G15").Value = 1.5
$ws.Range("G16").Value = 1.5
$ws.Range("G17").Value = 1
$ws.Range("G18").Value = 1
$ws.Range("G21").Value = 1.5

# --- Update STATUS (column E) for row 19 ---
$ws.Range("E19").Value = "Ready to Start"

# --- Move the active selection / visible top-left cell on the sheet ---
$ws.Activate() | Out-Null
$ws.Range("G19").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 7 | Out-Null
